$wb = $excel.ActiveWorkbook

# --- "Ranges" sheet: update Toronto Island and KW minimum date/gust values ---
$wsRanges = $wb.Worksheets.Item("Ranges")
$wsRanges.Range("A4").Value = 43117   # Toronto Island - Min_Date_Time
$wsRanges.Range("B4").Value = 67      # Toronto Island - Min_SpdOfMaxGust_km_h_
$wsRanges.Range("A7").Value = 37548   # KW - Min_Date_Time
$wsRanges.Range("B7").Value = 63      # KW - Min_SpdOfMaxGust_km_h_

# --- "Minimum" sheet: same updates, then re-sort the table ascending by Min speed ---
$wsMin = $wb.Worksheets.Item("Minimum")
$wsMin.Range("A2").Value = 43117      # Toronto Island - Min_Date_Time
$wsMin.Range("B2").Value = 67         # Toronto Island - Min_SpdOfMaxGust_km_h_
$wsMin.Range("A3").Value = 37548      # KW - Min_Date_Time
$wsMin.Range("B3").Value = 63         # KW - Min_SpdOfMaxGust_km_h_

$dataRange = $wsMin.Range("A2:C9")
$dataRange.Sort($wsMin.Range("B2:B9"), 1)
